$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update re-worded / corrected ticket descriptions (shared-string text edits) ---
# (Order matters for how the shared-string table gets rebuilt on save, so these
# follow the same chronological order as the original authoring session.)

# G30: HARDWARE note - 20*fsignal/24kHz -> 10*fsignal/12kHz
$ws.Range("G30").Value = "HARDWARE - Filtrer le signal provenant du compteur, en aval de l'opto. Simu à faire sous PSIM (fc > 10*fsignal = 12kHz)"

# --- Fill in the new ticket row (row 31), which was previously a blank template row ---

$ws.Range("A31").Value = 26
$ws.Range("B31").Value = "CORRECTION"
$ws.Range("C31").Value = "MAJEUR"
$ws.Range("D31").Value = "-"
$ws.Range("E31").Value = "X"
$ws.Range("F31").Value = "Ouvert"
$ws.Range("G31").Value = "Remplacer le - de Papp dans HCHP par 0`nRetirer les lignes vides à la fin lors d'un relevé HCHP (dernière ligne + celle avant mot d'état/Imax)"

# H31 needs the same date number-format as the sibling rows (e.g. H30) before
# writing the date value, otherwise it keeps the row's default "General" style.
$ws.Range("H30").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("H31").Value = 41883

$ws.Rows(31).RowHeight = 45

# G27: "Contraire" -> "Contraindre" (typo fix), text otherwise unchanged
$ws.Range("G27").Value = "Contraindre le changement de jour pour garantir la cloture du fichier`nEnvisager de tester 00:01:00 ?!"

# G28: add trailing clarification "de l'heure entière"
$ws.Range("G28").Value = "Enregistrer les index à +/- 1 seconde de l'heure entière"

# G29: "n'ai pas" -> "n'est pas" (typo fix)
$ws.Range("G29").Value = "Lors de l'enregistrement du courant max; veiller à ce qu'il n'est pas déjà été remis à 0"

# --- View state: scrolled/frozen pane position & active selection moved down ---
$ws.Range("B24").Select()
$ws.Range("B32").Select()

$wb.Save()
